$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("colName", "shortName", "janName", "statId", "group", "header"),
    @("pass_yds", "PASSING YDS", "passing_yds", 3, "Passing", "YDS"),
    @("pass_tds", "PASSING TD", "passing_td", 4, "Passing", "TD"),
    @("pass_int", "PASSING INT", "passing_int", 20, "Passing", "INT"),
    @("pass_att", "PASSING ATT", "passing_att", 1, "Passing", "A"),
    @("pass_comp", "PASSING COMP", "passing_comp", 0, "Passing", "C"),
    @("rush_att", "RUSHING RUSH", "rushing_rush", 23, "rushing", "CAR"),
    @("rush_yds", "RUSHING YDS", "rushing_yds", 24, "rushing", "YDS"),
    @("rush_tds", "RUSHING TD", "rushing_td", 25, "rushing", "TD"),
    @("rec", "RECEIVING REC", "receiving_rec", 53, "rec", "REC"),
    @("rec_yds", "RECEIVING YDS", "receiving_yds", 42, "rec", "YDS"),
    @("rec_tds", "RECEIVING TD", "receiving_td", 43, "rec", "TD"),
    @("dst_tackle", "DEFENSIVE TT", "defensive_tt", 109, "Defensive", "TT"),
    @("dst_sacks", "DEFENSIVE SCK", "defensive_sck", 99, "Defensive", "SCK"),
    @("dst_fum_force", "DEFENSIVE FF", "defensive_ff", 106, "Defensive", "FF"),
    @("dst_fum_rec", "DEFENSIVE FR", "defensive_fr", 96, "Defensive", "FR"),
    @("dst_int", "DEFENSIVE INT", "defensive_int", 95, "Defensive", "INT"),
    @("dst_int_tds", "TD RETURNS ITD", "td_returns_itd", 103, "TD Return", "ITD"),
    @("dst_ret_tds", "TD RETURNS FTD", "td_returns_ftd", 104, "TD Return", "FTD"),
    @("site_pts", "TOTAL PTS", "total_pts", "appliedTotal", "TOTAL", "FPTS"),
    @("fg", "FG FGM", "fg_fgm", 80, "kicking", "FG39"),
    @("fg_att", "FG FGA", "fg_fga", 81, "kicking", "FGA39"),
    @("fg_4049", "FG 40-49", "fg_40_49", 77, "kicking", "FG49"),
    @("fg_50", "FG 50", "fg_50", 74, "kicking", "FG50+"),
    @("fg_att_4049", "FG 40-49A", "fg_40_49a", 78, "kicking", "FGA49"),
    @("fg_att_50", "FG 50+A", "fg_50_a", 75, "kicking", "FGA50+"),
    @("xp", "XP XPM", "xp_xpm", 86, "kicking", "XP"),
    @("xp_att", "XP XPA", "xp_xpa", 87, "kicking", "XPA"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}

$null = $ws.Range("D13").Select()
